# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Nuevos periodos de mora (columna E), en orden ascendente, filas 16-23
$periodos = @("1910","1911","1912","2001","2002","2003","2004","2005")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Nuevo valor mora (columna F): todas las filas quedan en 34000 excepto la
# ultima (fila 23, periodo 2005) que conserva el valor 21533
for ($row = 16; $row -le 22; $row++) {
    $ws.Range("F$row").Value = 34000
}
$ws.Range("F23").Value = 21533

# Nuevo salario basico (columna G) para todas las filas de la tabla
for ($row = 16; $row -le 23; $row++) {
    $ws.Range("G$row").Value = 850000
}
